# Slide 3 ("STGraph - Conceptualization"), Content Placeholder 2:
# Prefix each of the two research-question paragraphs with a bold
# "RQ1" / "RQ2" label followed by " - " before the existing sentence.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$tr = $s.Shapes.Item(2).TextFrame.TextRange

# --- Paragraph 1: "Can we separate data by temporal granularity ..." ---
$para1 = $tr.Paragraphs(1, 1)
$body1 = $tr.Characters($para1.Start, 118)
$body1.Text = " - " + $body1.Text
$para1 = $tr.Paragraphs(1, 1)
$para1.InsertBefore("RQ1") | Out-Null
$label1 = $tr.Characters($para1.Start, 3)
$label1.Font.Bold = $true

# --- Paragraph 2: "If so, can we embed two different data-layout ..." ---
$para2 = $tr.Paragraphs(2, 1)
$body2 = $tr.Characters($para2.Start, 118)
$body2.Text = " - " + $body2.Text
$para2 = $tr.Paragraphs(2, 1)
$para2.InsertBefore("RQ2") | Out-Null
$label2 = $tr.Characters($para2.Start, 3)
$label2.Font.Bold = $true
